$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R1")  # "Active_Outages" / sheet1.xml region R1

# Append a new row (row 6) of data to the Active Outages sheet, mirroring
# the existing rows' layout: only Region (B), Hub Site (D), Power Source
# (I), Battery Backup Status (J) and Site Owner (L) are populated; the
# remaining columns (A, C, E, F, G, H, K) stay blank, same as rows 4-5.
$ws.Cells.Item(6, 1).Value = ""
$ws.Cells.Item(6, 2).Value = "R4"
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = "JED0123"
$ws.Cells.Item(6, 5).Value = ""
$ws.Cells.Item(6, 6).Value = ""
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).Value = ""
$ws.Cells.Item(6, 9).Value = "SCECO"
$ws.Cells.Item(6, 10).Value = "In progress"
$ws.Cells.Item(6, 11).Value = ""
$ws.Cells.Item(6, 12).Value = "Latis"
